$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$d = $p.Designs.Item(1)
try {
  $s.Design = $d
  Write-Output "OK"
} catch {
  Write-Output "ERROR: $($_.Exception.Message)"
}
